$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1452
$ws.Range("J40").Value = 1452
$ws.Range("L40").Value = 1452
$ws.Range("N40").Value = -1802
$ws.Range("H74").Value = 15211.143
$ws.Range("J74").Value = 29799.8
$ws.Range("L74").Value = 29799.8
$ws.Range("N74").Value = -31671.8
$ws.Range("H77").Value = 15211.143
$ws.Range("J77").Value = 29799.8
$ws.Range("L77").Value = 148999
$ws.Range("N77").Value = -158359
$ws.Range("H98").Value = 1013.7857
$ws.Range("I98").Value = 982.4231
$ws.Range("J98").Value = 1421.5
$ws.Range("K98").Value = 982.4231
$ws.Range("L98").Value = 1421.5
$ws.Range("M98").Value = 515.5769
$ws.Range("N98").Value = -4417.5
$ws.Range("H122").Value = 1013.7857
$ws.Range("I122").Value = 982.4231
$ws.Range("J122").Value = 1421.5
$ws.Range("K122").Value = 2947.2693
$ws.Range("L122").Value = 4264.5
$ws.Range("M122").Value = -497.2692999999999
$ws.Range("N122").Value = -9164.5
$ws.Range("H132").Value = 9525450
$ws.Range("I132").Value = 10102356
$ws.Range("K132").Value = 30307068
$ws.Range("M132").Value = -30304538
$ws.Range("H141").Value = 1736.25
$ws.Range("I141").Value = 1594.3478
$ws.Range("K141").Value = 4783.0434
$ws.Range("M141").Value = 396.9565999999995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 35716612
$ws.Range("I2").Value = 41668828
$ws.Range("K2").Value = 41668828
$ws.Range("M2").Value = -41668715
$ws.Range("H14").Value = 9999.666999999999
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 9999.666999999999
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 9999.666999999999
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10349.667
$ws.Range("H32").Value = 4104.5083
$ws.Range("I32").Value = 3083.566
$ws.Range("K32").Value = 3083.566
$ws.Range("M32").Value = -2796.566
$ws.Range("H36").Value = 26940.875
$ws.Range("I36").Value = 2921.3333
$ws.Range("K36").Value = 2921.3333
$ws.Range("M36").Value = -2575.3333
$ws.Range("H110").Value = 7959.476
$ws.Range("I110").Value = 11468.091
$ws.Range("J110").Value = 4100
$ws.Range("K110").Value = 11468.091
$ws.Range("L110").Value = 4100
$ws.Range("M110").Value = -9423.091
$ws.Range("N110").Value = -8190
$ws.Range("H116").Value = 35716612
$ws.Range("I116").Value = 41668828
$ws.Range("K116").Value = 41668828
$ws.Range("M116").Value = -41666534
$ws.Range("H122").Value = 1762.9048
$ws.Range("I122").Value = 1481.875
$ws.Range("K122").Value = 4445.625
$ws.Range("M122").Value = -1995.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 35716612
$ws.Range("I3").Value = 41668828
$ws.Range("K3").Value = 41668828
$ws.Range("M3").Value = -41668714
$ws.Range("H105").Value = 1808
$ws.Range("I105").Value = 1837.1666
$ws.Range("K105").Value = 1837.1666
$ws.Range("M105").Value = -90.16660000000002
$ws.Range("H134").Value = 1504.1372
$ws.Range("I134").Value = 1555.4565
$ws.Range("K134").Value = 4666.3695
$ws.Range("M134").Value = -2131.3695

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 20000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H113").Value = 20000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 1233.9375
$ws.Range("I122").Value = 1247.2727
$ws.Range("J122").Value = 1204.6
$ws.Range("K122").Value = 3741.8181
$ws.Range("L122").Value = 3613.8
$ws.Range("M122").Value = -1291.8181
$ws.Range("N122").Value = -8513.799999999999
$ws.Range("H132").Value = 4215.8
$ws.Range("I132").Value = 4180.5557
$ws.Range("K132").Value = 12541.6671
$ws.Range("M132").Value = -10011.6671

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.90625
$ws.Range("I2").Value = 60.47826
$ws.Range("J2").Value = 15.777778
$ws.Range("K2").Value = 362.86956
$ws.Range("L2").Value = 94.666668
$ws.Range("M2").Value = -249.86956
$ws.Range("N2").Value = -320.666668
$ws.Range("H5").Value = 28433.334
$ws.Range("I5").Value = 300
$ws.Range("J5").Value = 42500
$ws.Range("K5").Value = 900
$ws.Range("L5").Value = 127500
$ws.Range("M5").Value = -788
$ws.Range("N5").Value = -127724
$ws.Range("H32").Value = 1950
$ws.Range("J32").Value = 1850
$ws.Range("L32").Value = 5550
$ws.Range("N32").Value = -6116
$ws.Range("H94").Value = 5157.5835
$ws.Range("I94").Value = 2799.3333
$ws.Range("K94").Value = 8397.999899999999
$ws.Range("M94").Value = -7721.999899999999
$ws.Range("H105").Value = 14000
$ws.Range("J105").Value = 14000
$ws.Range("L105").Value = 42000
$ws.Range("N105").Value = -47242
$ws.Range("H135").Value = 28433.334
$ws.Range("I135").Value = 300
$ws.Range("J135").Value = 42500
$ws.Range("K135").Value = 2700
$ws.Range("L135").Value = 382500
$ws.Range("M135").Value = -165
$ws.Range("N135").Value = -387570

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 50000
$ws.Range("J15").Value = 50000
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50576
$ws.Range("H19").Value = 4669
$ws.Range("I19").Value = 4889.5454
$ws.Range("J19").Value = 4062.5
$ws.Range("K19").Value = 4889.5454
$ws.Range("L19").Value = 4062.5
$ws.Range("M19").Value = -4601.5454
$ws.Range("N19").Value = -4638.5
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H126").Value = 22279.412
$ws.Range("J126").Value = 4176.625
$ws.Range("L126").Value = 12529.875
$ws.Range("N126").Value = -17469.875
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4831.25
$ws.Range("I122").Value = 4081.8823
$ws.Range("J122").Value = 5989.364
$ws.Range("K122").Value = 12245.6469
$ws.Range("L122").Value = 17968.092
$ws.Range("M122").Value = -9795.6469
$ws.Range("N122").Value = -22868.092

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 8329.333000000001
$ws.Range("I39").Value = 9994
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 9994
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = -9581
$ws.Range("N39").Value = -5826
